$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.872.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.229.18'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.49%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.64'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.90%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.63'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +18.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0975'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '57.97'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +7.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.561.86'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.96'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.866'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.226.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.834.96'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0974'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.11'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.90'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +11.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.98%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.54'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +8.02%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.83'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.89'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.58'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0730'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.70'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.04'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +24.07%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +10.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0298'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +12.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.29'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.82%  '
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.99'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '68.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.06'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +22.56%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +11.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.93'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.77'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.70'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.13%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +7.74%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.59%  '
